$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($row, $vals) {
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, 4 + $i).Value = $vals[$i]
    }
}

Set-RowValues 8 @("فصل سوم منتهی به 1399/09", "فصل چهارم منتهی به 1399/12", "فصل اول منتهی به 1400/03", "فصل دوم منتهی به 1400/06", "فصل سوم منتهی به 1400/09", "فصل چهارم منتهی به 1400/12", "فصل اول منتهی به 1401/03", "فصل دوم منتهی به 1401/06", "فصل سوم منتهی به 1401/09", "فصل چهارم منتهی به 1401/12")
Set-RowValues 9 @("1400-11-02 (2)", "1401-02-25 (12)", "1401-05-01 (3)", "1401-08-30 (4)", "1401-11-19 (3)", "1402-02-27 (12)", "1401-05-01", "1401-08-30 (2)", "1401-11-19 (2)", "1402-02-27 (3)")
Set-RowValues 12 @(1213811, 1099830, 1513024, 1717640, -333461, 2436389, 2683124, 1739073, 1151190, 3459418)
Set-RowValues 13 @(-192397, 61075, 0, -178414, -347236, 136680, 0, -279401, -158450, -248155)
Set-RowValues 14 @(1021414, 1160905, 1513024, 1539226, -680697, 2573069, 2683124, 1459672, 992740, 3211263)
Set-RowValues 16 @(3333, 1700, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 17 @(-19581, -28411, -10576, -16051, 0, -25455, 0, -20555, -522, -222051)
Set-RowValues 18 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 19 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 20 @(0, 0, 0, 0, 0, 0, 0, 0, 0, -3426)
Set-RowValues 21 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 22 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 23 @(0, 0, 0, 615600, 29535, -645135, 0, 0, 0, 0)
Set-RowValues 24 @(0, 0, 0, 0, -76020, 76020, 0, 0, 0, 0)
Set-RowValues 25 @(0, 382747, 0, 0, 0, 11851, 0, 0, 0, 0)
Set-RowValues 26 @(-1150198, -1080051, -2130113, 2054093, 76020, -1687998, -746398, 199992, -44970, -3965550)
Set-RowValues 27 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 28 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 29 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 30 @(-578, 0, 0, 24844, 0, 0, 89971, -48603, 48360, 15757)
Set-RowValues 31 @(30292, 54526, 25444, 123012, 166425, 93549, 158484, 268494, 195174, 294382)
Set-RowValues 32 @(-1136732, -669489, -2115245, 2801498, 195960, -2177168, -497943, 399328, 198042, -3880888)
Set-RowValues 33 @(-115318, 491416, -602221, 4340724, -484737, 395901, 2185181, 1859000, 1190782, -669625)
Set-RowValues 35 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 36 @("-", "-", "-", "-", 0, 0, "-", "-", 0, 0)
Set-RowValues 37 @(0, 224434, 0, 92772, -92772, 149376, 0, 80512, 0, 68448)
Set-RowValues 38 @(0, -290605, 0, -63237, 63237, -147929, 0, -14512, 0, -75546)
Set-RowValues 39 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 40 @(0, -17160, 0, -17160, 0, -17160, 0, -17160, 0, -17160)
Set-RowValues 41 @(321, -21748, 0, -23046, 0, -24323, 0, -25620, 0, -26897)
Set-RowValues 42 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 43 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 44 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 45 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 46 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 47 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 48 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 49 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues 50 @(-94, -27006, 0, -3131441, -24148, -341, 0, -1999558, -2856701, -32355)
Set-RowValues 51 @(227, -132085, 0, -3142112, -53683, -40377, 0, -1976338, -2856701, -83510)
Set-RowValues 52 @(-115091, 359331, -602221, 1198612, -538420, 355524, 2185181, -117338, -1665919, -753135)
Set-RowValues 53 @(628737, 539172, 819625, 210560, 1421181, 882761, 1224333, 3409514, 3292993, 1708760)
Set-RowValues 54 @(25526, -78878, -6844, 12009, 0, -13952, 0, 817, 81686, 296209)
Set-RowValues 55 @(539172, 819625, 210560, 1421181, 882761, 1224333, 3409514, 3292993, 1708760, 1251834)
Set-RowValues 56 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)

# Column widths: the quarter table shifted one column to the left (oldest
# quarter dropped, newest quarter appended), so the "wide" (Q4-style, 31 chars)
# columns move from F/J to E/I, plus the brand-new last column M also becomes
# a Q4 column and gets the wide width.
$ws.Columns("D").ColumnWidth = 29 - (5/6)
$ws.Columns("E").ColumnWidth = 31 - (5/6)
$ws.Columns("F").ColumnWidth = 29 - (5/6)
$ws.Columns("G").ColumnWidth = 29 - (5/6)
$ws.Columns("H").ColumnWidth = 29 - (5/6)
$ws.Columns("I").ColumnWidth = 31 - (5/6)
$ws.Columns("J").ColumnWidth = 29 - (5/6)
$ws.Columns("K").ColumnWidth = 29 - (5/6)
$ws.Columns("L").ColumnWidth = 29 - (5/6)
$ws.Columns("M").ColumnWidth = 31 - (5/6)
